$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "87.514.32"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").Value = "3.165.93"
$ws.Range("E3").Value = "  -6.52%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.13"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -5.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "608.07"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -5.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.375"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -10.82%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.663"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").Value = "3.161.14"
$ws.Range("E10").Value = "  -6.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.533"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -14.06%  "
$ws.Range("E12").Value = "  +4.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000239"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -16.24%  "
$ws.Range("D14").Value = "3.730.86"
$ws.Range("E14").Value = "  -7.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.22"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -6.01%  "
$ws.Range("D16").Value = "87.177.42"
$ws.Range("E16").Value = "  -1.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "32.02"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -12.12%  "
$ws.Range("D18").Value = "3.166.49"
$ws.Range("E18").Value = "  -6.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.05"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.35"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -10.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "413.23"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -9.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.42"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -12.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.03"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -10.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.14"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -6.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.78"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -6.42%  "
$ws.Range("D26").Value = "3.324.40"
$ws.Range("E26").Value = "  -7.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "73.00"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -9.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000129"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -9.94%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  -0.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.158"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -16.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "543.00"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.17"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -13.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.29"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -16.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.84"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -12.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.66"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -9.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.132"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -7.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "21.71"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -9.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "21.81"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.97"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.58%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.90"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -9.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.366"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -15.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "148.76"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -5.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "171.85"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -8.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.13"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -7.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.124"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.21"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -14.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.94"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -12.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.691"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -12.17%  "
